# Reapply "new changes in ops (ordercreation & orderpage & order form)"
# Adds Typist / Typist QC / Lob / Process / Product Name / State / County /
# Municipality / Status / Tier columns (E:O) to the existing order row, and
# appends a second order row (row 3, a "Typing" order) below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1): extend from D1 out to O1, copying the existing
#    header format (bold font + fill + border) from D1.
# ---------------------------------------------------------------------
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:O1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"
$ws.Range("G1").Value = "Client"
$ws.Range("H1").Value = "Lob"
$ws.Range("I1").Value = "Process"
$ws.Range("J1").Value = "Product Name"
$ws.Range("K1").Value = "State"
$ws.Range("L1").Value = "County"
$ws.Range("M1").Value = "Municipality"
$ws.Range("N1").Value = "Status"
$ws.Range("O1").Value = "Tier"

# ---------------------------------------------------------------------
# 2. Row 2 (existing "Search" order): the old E2 ("Client"/Beeline) slides
#    right to G2; columns E2/F2 become blank Typist/Typist QC cells; and
#    H2:O2 are filled in with the Lob/Process/.../Tier data that used to
#    live in F2:M2.
# ---------------------------------------------------------------------

# Donor cells (already present in the sheet) to copy formats from:
#   s=1 -> B2, s=4 -> C2, s=5 -> D2
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2:F2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("H2:O2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# Clear leftover text so the freshly-formatted E2/F2 stay blank.
$ws.Range("E2").ClearContents() | Out-Null
$ws.Range("F2").ClearContents() | Out-Null

$ws.Range("G2").Value = "Beeline"
$ws.Range("H2").Value = "Title"
$ws.Range("I2").Value = "Search"
$ws.Range("J2").Value = "Current Owner Search"
$ws.Range("K2").Value = "AL"
$ws.Range("L2").Value = "Shelby"
$ws.Range("M2").Value = "ALShelby"
$ws.Range("N2").Value = "WIP"
$ws.Range("O2").Value = "Search(T1)"

# ---------------------------------------------------------------------
# 3. Row 3 (new "Typing" order) - build the whole row from the row-2
#    formats, then fill in its values.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3:F3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("G2").Copy() | Out-Null
$ws.Range("G3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("H2:O2").Copy() | Out-Null
$ws.Range("H3:O3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# E3 starts from the full thin-box border (same as H2's style) and then
# drops its top edge (it sits directly under row 2's bottom border).
$ws.Range("H2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("E3").Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone

$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("D3").ClearContents() | Out-Null

$ws.Range("A3").Value = 45439.083333333336
$ws.Range("B3").Value = "Be18-002"
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL5317"
$ws.Range("G3").Value = "Beeline"
$ws.Range("H3").Value = "Title"
$ws.Range("I3").Value = "Typing"
$ws.Range("J3").Value = "Commitment Typing"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "FLClay"
$ws.Range("N3").Value = "Typing"
$ws.Range("O3").Value = "Typing(T1)"

# ---------------------------------------------------------------------
# 4. Column widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.6666666666667    # C: 36.33 -> ~20.55
$ws.Columns.Item(6).ColumnWidth = 15.1666666666667    # F: 16 (unchanged)
$ws.Columns.Item(7).ColumnWidth = 15.1666666666667    # G: 16 (unchanged)
$ws.Columns.Item(8).ColumnWidth = 15.1666666666667    # H: 29.55 -> 16
$ws.Columns.Item(10).ColumnWidth = 18.3333333333333   # J: new -> ~19.11

# ---------------------------------------------------------------------
# 5. Selection, to match the saved view.
# ---------------------------------------------------------------------
$ws.Range("J6").Select() | Out-Null
